$wb = $excel.ActiveWorkbook

$wsChart = $wb.Worksheets("Chart")
$wsTable = $wb.Worksheets("Table")

# --- Sheet "Chart" (sheet1.xml): append 4 new daily rows (28-31) ---

# Column A holds date strings stored as TEXT (shared strings), not real Excel
# dates, so force a text number format before assigning the date-like
# strings - otherwise they'd be auto-parsed into date serial numbers.
$wsChart.Range("A28:A31").NumberFormat = "@"

$wsChart.Range("A28").Value = "2025-10-31"
$wsChart.Range("B28").Value = 23
$wsChart.Range("C28").Value = 0
$wsChart.Range("D28").Value = 0

$wsChart.Range("A29").Value = "2025-11-01"
$wsChart.Range("B29").Value = 23
$wsChart.Range("C29").Value = 0
$wsChart.Range("D29").Value = 0

$wsChart.Range("A30").Value = "2025-11-02"
$wsChart.Range("B30").Value = 23
$wsChart.Range("C30").Value = 0
$wsChart.Range("D30").Value = 0

$wsChart.Range("A31").Value = "2025-11-03"
$wsChart.Range("B31").Value = 23
$wsChart.Range("C31").Value = 0
# D31 has no Impressions figure yet for the newest day (exported as an empty
# text cell, matching the blank placeholder cells already used elsewhere in
# this sheet, e.g. B2/C2). A leading apostrophe forces an empty *text* cell
# instead of clearing the cell outright.
$wsChart.Range("D31").Value = "'"

# The quick NumberFormat/apostrophe tricks above leave behind "number/text
# stored as text" styling on the affected cells; copy the plain (default)
# formatting back from already-existing neighbour cells so the new cells
# match the rest of the sheet's styling.
$wsChart.Range("A27").Copy()
$wsChart.Range("A28:A31").PasteSpecial(-4122)
$wsChart.Range("D27").Copy()
$wsChart.Range("D31").PasteSpecial(-4122)

# --- Sheet "Table" (sheet2.xml): refreshed Videos count (22 -> 23) ---
$wsTable.Range("C2").Value = 23
